$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row): extend sequence with P1=14, Q1=15, copying O1's style ---
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Rows 2-25: swap I<->K and M<->O values, then add P and Q columns with value 2 ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I: was 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: was 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: was 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: was 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P: new column
    $ws.Cells.Item($r, 17).Value = 2  # Q: new column
}
